# Adding open account test case
#
# 1) addCustomer sheet: drop the unused "alerttext" column (D) and add
#    three new customer rows that already have a currency assigned.
# 2) Add a brand-new "OpenAccount" sheet with its own small data set.

$wb = $excel.ActiveWorkbook
$ws1 = $wb.ActiveSheet   # addCustomer

# --- Remove the alerttext column (D) from addCustomer -----------------
$ws1.Range("D1").EntireColumn.Delete()

# --- Append the new rows to addCustomer --------------------------------
$ws1.Range("A3").Value = "Hamon"
$ws1.Range("B3").Value = "Roger"
$ws1.Range("C3").Value = 652345

$ws1.Range("A4").Value = "Soraya"
$ws1.Range("B4").Value = "Gregor"
$ws1.Range("C4").Value = 6646753

$ws1.Range("A5").Value = "Khaloy"
$ws1.Range("B5").Value = "Rodrigues"
$ws1.Range("C5").Value = 3425131

# --- Add the new OpenAccount sheet after addCustomer --------------------
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "OpenAccount"

$ws2.Range("A1").Value = "Customer"
$ws2.Range("B1").Value = "Currency"
$ws2.Range("A2").Value = "Hermoine Granger"
$ws2.Range("B2").Value = "Rupee"

# --- Restore addCustomer as the active sheet, cursor left on C4 --------
$ws1.Activate()
$ws1.Range("C4").Select()
